$wb = $excel.ActiveWorkbook

# Update values on the "additional" sheet (4th sheet) to match PDF edition 10
$ws = $wb.Worksheets.Item("additional")
$ws.Range("B2").Value = 47.2
$ws.Range("B3").Value = 82.7
$ws.Range("B4").Value = 61.6

# Make "additional" the active sheet (tabSelected / activeTab in the OOXML)
# and move its selection to B5.
$ws.Activate()
$ws.Range("B5").Select()
